$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update question text in B2 (was "How many people can I hire?")
$ws.Range("B2").Value = "I have a client who looking to hire a Chinese national with a recent masters in engineering. Is he eligible for support?"

# B3 keeps the same visible text "Is the grant available to a startup?"
$ws.Range("B3").Value = "Is the grant available to a startup?"

# Update selected cell from C7 to B3
$ws.Range("B3").Select()
